$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1255.2174
$ws.Range("I17").Value = 664.44446
$ws.Range("J17").Value = 1635
$ws.Range("K17").Value = 1993.33338
$ws.Range("L17").Value = 4905
$ws.Range("M17").Value = -1825.33338
$ws.Range("N17").Value = -5241

$ws.Range("H33").Value = 4088.3845
$ws.Range("I33").Value = 81.625
$ws.Range("J33").Value = 10499.2
$ws.Range("K33").Value = 81.625
$ws.Range("L33").Value = 10499.2
$ws.Range("M33").Value = 147.375
$ws.Range("N33").Value = -10957.2

$ws.Range("H62").Value = 8920.5625
$ws.Range("I62").Value = 11181.818
$ws.Range("J62").Value = 3945.8
$ws.Range("K62").Value = 11181.818
$ws.Range("L62").Value = 3945.8
$ws.Range("M62").Value = -10557.818
$ws.Range("N62").Value = -5193.8

$ws.Range("H65").Value = 8920.5625
$ws.Range("I65").Value = 11181.818
$ws.Range("J65").Value = 3945.8
$ws.Range("K65").Value = 55909.09
$ws.Range("L65").Value = 19729
$ws.Range("M65").Value = -52789.09
$ws.Range("N65").Value = -25969

$ws.Range("H98").Value = 4119992.8
$ws.Range("I98").Value = 6124.722
$ws.Range("J98").Value = 12347729
$ws.Range("K98").Value = 6124.722
$ws.Range("L98").Value = 12347729
$ws.Range("M98").Value = -4626.722
$ws.Range("N98").Value = -12350725

$ws.Range("H100").Value = 2525.3684
$ws.Range("I100").Value = 1342.2307
$ws.Range("K100").Value = 1342.2307
$ws.Range("M100").Value = -801.2307000000001

$ws.Range("H103").Value = 71429300
$ws.Range("I103").Value = 876.8
$ws.Range("J103").Value = 250000350
$ws.Range("K103").Value = 2630.4
$ws.Range("L103").Value = 750001050
$ws.Range("M103").Value = -2044.4
$ws.Range("N103").Value = -750002222

$ws.Range("H106").Value = 3889.4443
$ws.Range("I106").Value = 3784.2354
$ws.Range("K106").Value = 3784.2354
$ws.Range("M106").Value = -3153.2354

$ws.Range("H122").Value = 4119992.8
$ws.Range("I122").Value = 6124.722
$ws.Range("J122").Value = 12347729
$ws.Range("K122").Value = 18374.166
$ws.Range("L122").Value = 37043187
$ws.Range("M122").Value = -15924.166
$ws.Range("N122").Value = -37048087

$ws.Range("H132").Value = 4020.1428
$ws.Range("I132").Value = 4254.5835
$ws.Range("J132").Value = 2613.5
$ws.Range("K132").Value = 12763.7505
$ws.Range("L132").Value = 7840.5
$ws.Range("M132").Value = -10233.7505
$ws.Range("N132").Value = -12900.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1788.6086
$ws.Range("I74").Value = 1233.8
$ws.Range("J74").Value = 2215.3845
$ws.Range("K74").Value = 1233.8
$ws.Range("L74").Value = 2215.3845
$ws.Range("M74").Value = -359.8
$ws.Range("N74").Value = -3963.3845

$ws.Range("H77").Value = 1788.6086
$ws.Range("I77").Value = 1233.8
$ws.Range("J77").Value = 2215.3845
$ws.Range("K77").Value = 6169
$ws.Range("L77").Value = 11076.9225
$ws.Range("M77").Value = -1801
$ws.Range("N77").Value = -19812.9225

$ws.Range("H102").Value = 4717.278
$ws.Range("I102").Value = 3825
$ws.Range("K102").Value = 3825
$ws.Range("M102").Value = -2203

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1546.2084
$ws.Range("I99").Value = 1135.3125
$ws.Range("K99").Value = 1135.3125
$ws.Range("M99").Value = 362.6875

$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2998
$ws.Range("I16").Value = 2747.5
$ws.Range("J16").Value = 4000
$ws.Range("K16").Value = 2747.5
$ws.Range("L16").Value = 4000
$ws.Range("M16").Value = -2460.5
$ws.Range("N16").Value = -4574

$ws.Range("H99").Value = 2343.8215
$ws.Range("I99").Value = 2312.652
$ws.Range("K99").Value = 2312.652
$ws.Range("M99").Value = -814.652

$ws.Range("H113").Value = 2998
$ws.Range("I113").Value = 2747.5
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 2747.5
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -577.5
$ws.Range("N113").Value = -8340

$ws.Range("H126").Value = 2343.8215
$ws.Range("I126").Value = 2312.652
$ws.Range("K126").Value = 6937.956
$ws.Range("M126").Value = -4467.956

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 743.08
$ws.Range("I131").Value = 330.69232
$ws.Range("J131").Value = 804.7012
$ws.Range("K131").Value = 992.07696
$ws.Range("L131").Value = 2414.1036
$ws.Range("M131").Value = 4047.92304
$ws.Range("N131").Value = -12494.1036

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H58").Value = 6418.7
$ws.Range("I58").Value = 1241
$ws.Range("J58").Value = 6994
$ws.Range("K58").Value = 1241
$ws.Range("L58").Value = 6994
$ws.Range("M58").Value = -964
$ws.Range("N58").Value = -7548

$ws.Range("H97").Value = 1436.6666
$ws.Range("I97").Value = 763.38464
$ws.Range("J97").Value = 3187.2
$ws.Range("K97").Value = 763.38464
$ws.Range("L97").Value = 3187.2
$ws.Range("M97").Value = -267.38464
$ws.Range("N97").Value = -4179.2

$ws.Range("H113").Value = 38462730
$ws.Range("I113").Value = 983.3077
$ws.Range("K113").Value = 983.3077
$ws.Range("M113").Value = 1186.6923

$ws.Range("H126").Value = 10418307
$ws.Range("I126").Value = 1889
$ws.Range("J126").Value = 20834726
$ws.Range("K126").Value = 5667
$ws.Range("L126").Value = 62504178
$ws.Range("M126").Value = -3197
$ws.Range("N126").Value = -62509118

$ws.Range("H132").Value = 3985.7678
$ws.Range("I132").Value = 4143.2383
$ws.Range("J132").Value = 3513.3572
$ws.Range("K132").Value = 12429.7149
$ws.Range("L132").Value = 10540.0716
$ws.Range("M132").Value = -9899.714899999999
$ws.Range("N132").Value = -15600.0716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1951.6666
$ws.Range("I61").Value = 1108.909
$ws.Range("J61").Value = 2878.7
$ws.Range("K61").Value = 1108.909
$ws.Range("L61").Value = 2878.7
$ws.Range("M61").Value = -906.9090000000001
$ws.Range("N61").Value = -3282.7

$ws.Range("H100").Value = 2761.238
$ws.Range("I100").Value = 1926.8182
$ws.Range("K100").Value = 1926.8182
$ws.Range("M100").Value = -1385.8182

$ws.Range("H113").Value = 1951.6666
$ws.Range("I113").Value = 1108.909
$ws.Range("J113").Value = 2878.7
$ws.Range("K113").Value = 1108.909
$ws.Range("L113").Value = 2878.7
$ws.Range("M113").Value = 1061.091
$ws.Range("N113").Value = -7218.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2774.2942
$ws.Range("J81").Value = 2878
$ws.Range("L81").Value = 5756
$ws.Range("N81").Value = -7878

$ws.Range("H84").Value = 2774.2942
$ws.Range("J84").Value = 2878
$ws.Range("L84").Value = 28780
$ws.Range("N84").Value = -39388

$ws.Range("H113").Value = 25641516
$ws.Range("I113").Value = 409.45456
$ws.Range("K113").Value = 1228.36368
$ws.Range("M113").Value = 941.6363200000001
